$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 379, pushing existing rows 379-396 down to 381-398.
$ws.Rows.Item(379).Insert()
$ws.Rows.Item(379).Insert()

# New row 379: "Pintón" quality, week of 2021-11-09 (serial 44509)
$ws.Range("A379").Value = 7
$ws.Range("B379").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C379").Value = "Ñuble"
$ws.Range("D379").Value = 44509
$ws.Range("E379").Value = 16
$ws.Range("F379").Value = "Fruta"
$ws.Range("G379").Value = 100108
$ws.Range("H379").Value = "Tropicales y subtropicales"
$ws.Range("I379").Value = 100108006
$ws.Range("J379").Value = "Plátano"
$ws.Range("K379").Value = "Sin especificar"
$ws.Range("L379").Value = "Pintón"
$ws.Range("M379").Value = 240
$ws.Range("N379").Value = 16000
$ws.Range("O379").Value = 17000
$ws.Range("P379").Value = 16500
$ws.Range("Q379").Value = "$/caja 20 kilos"
$ws.Range("R379").Value = "Ecuador"
$ws.Range("S379").Value = 825
$ws.Range("T379").Value = 20

# New row 380: "Primera Pintón" quality, same week (serial 44509)
$ws.Range("A380").Value = 7
$ws.Range("B380").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C380").Value = "Ñuble"
$ws.Range("D380").Value = 44509
$ws.Range("E380").Value = 16
$ws.Range("F380").Value = "Fruta"
$ws.Range("G380").Value = 100108
$ws.Range("H380").Value = "Tropicales y subtropicales"
$ws.Range("I380").Value = 100108006
$ws.Range("J380").Value = "Plátano"
$ws.Range("K380").Value = "Sin especificar"
$ws.Range("L380").Value = "Primera Pintón"
$ws.Range("M380").Value = 240
$ws.Range("N380").Value = 18000
$ws.Range("O380").Value = 19000
$ws.Range("P380").Value = 18500
$ws.Range("Q380").Value = "$/caja 20 kilos"
$ws.Range("R380").Value = "Ecuador"
$ws.Range("S380").Value = 925
$ws.Range("T380").Value = 20
